$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.838.16"
$ws.Range("E2").Value = "  +6.41%  "
$ws.Range("D3").Value = "2.056.56"
$ws.Range("E3").Value = "  +3.55%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "253.05"
$ws.Range("E5").Value = "  +4.53%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.654"
$ws.Range("E6").Value = "  +2.61%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "65.72"
$ws.Range("E7").Value = "  +15.88%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +5.84%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "59.99"
$ws.Range("E10").Value = "  -0.23%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0766"
$ws.Range("E11").Value = "  +5.13%  "
$ws.Range("E12").Value = "  +1.08%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.932"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.24"
$ws.Range("E14").Value = "  +8.19%  "
$ws.Range("D15").Value = "2.351.83"
$ws.Range("E15").Value = "  +3.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.60"
$ws.Range("E16").Value = "  +6.83%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "20.58"
$ws.Range("E17").Value = "  +20.50%  "
$ws.Range("D18").Value = "2.048.56"
$ws.Range("E18").Value = "  +3.12%  "
$ws.Range("D19").Value = "37.689.72"
$ws.Range("E19").Value = "  +6.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "73.89"
$ws.Range("E20").Value = "  +5.19%  "
$ws.Range("D21").Value = "0.0₃0884"
$ws.Range("E21").Value = "  +5.78%  "
$ws.Range("E22").Value = "  +6.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "238.47"
$ws.Range("E23").Value = "  +2.54%  "
$ws.Range("E24").Value = "  +15.44%  "
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.40"
$ws.Range("E26").Value = "  +3.90%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.65"
$ws.Range("E27").Value = "  +6.22%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "160.17"
$ws.Range("E28").Value = "  -2.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.99"
$ws.Range("E29").Value = "  +2.59%  "
$ws.Range("E30").Value = "  +9.46%  "
$ws.Range("E31").Value = "  +2.93%  "
$ws.Range("E32").Value = "  +24.76%  "
$ws.Range("E33").Value = "  +6.74%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.78"
$ws.Range("E34").Value = "  +12.00%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0618"
$ws.Range("E35").Value = "  +5.45%  "
$ws.Range("E36").Value = "  +4.16%  "
$ws.Range("E37").Value = "  +3.89%  "
$ws.Range("E38").Value = "  -0.10%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.00"
$ws.Range("E39").Value = "  +22.71%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.102"
$ws.Range("E40").Value = "  +14.69%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.83"
$ws.Range("E41").Value = "  +25.54%  "
$ws.Range("E42").Value = "  +4.61%  "
$ws.Range("E43").Value = "  +5.36%  "
$ws.Range("E44").Value = "  +3.99%  "
$ws.Range("E45").Value = "  +5.86%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.09"
$ws.Range("E46").Value = "  +10.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "17.03"
$ws.Range("E47").Value = "  +10.35%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "95.81"
$ws.Range("E48").Value = "  +5.47%  "
$ws.Range("D49").Value = "1.431.11"
$ws.Range("E49").Value = "  +4.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.95"
$ws.Range("E50").Value = "  +2.57%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "47.60"
$ws.Range("E51").Value = "  +4.29%  "
